# Update the Correspond Handoff/Handback datetimes for the
# f34f748b-... rows (row 5) on the "zh-cn" and "de-de" sheets, to
# reflect the regenerated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-19 07:19:04"
$wsZhCn.Range("G5").Value = "2016-01-19 07:19:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-19 07:19:14"
$wsDeDe.Range("G5").Value = "2016-01-19 07:20:04"
